# Fruta / hortaliza, semanal
# Two new weekly price records (Región del Maule, week of 2022-02-18 / serial 44610)
# are inserted at the top of the data block (rows 459-460), pushing every
# existing record below down by two rows (old row 461 -> 463, ..., old row 485 -> 487).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 459, shifting the rest down.
$ws.Rows("459:460").Insert()

# New row 459: Especial quality entry for the new date.
$ws.Range("A459").Value = 5
$ws.Range("B459").Value = "Macroferia Regional de Talca"
$ws.Range("C459").Value = "Maule"
$ws.Range("D459").Value = 44610
$ws.Range("E459").Value = 7
$ws.Range("F459").Value = "Fruta"
$ws.Range("G459").Value = 100101
$ws.Range("H459").Value = "Berries"
$ws.Range("I459").Value = 100112025
$ws.Range("J459").Value = "Frutilla"
$ws.Range("K459").Value = "Sin especificar"
$ws.Range("L459").Value = "Especial"
$ws.Range("M459").Value = 300
$ws.Range("N459").Value = 8000
$ws.Range("O459").Value = 8000
$ws.Range("P459").Value = 8000
$ws.Range("Q459").Value = "`$/bandeja 7 kilos"
$ws.Range("R459").Value = "Región del Maule"
$ws.Range("S459").Value = 1143
$ws.Range("T459").Value = 7

# New row 460: Primera quality entry for the same new date.
$ws.Range("A460").Value = 5
$ws.Range("B460").Value = "Macroferia Regional de Talca"
$ws.Range("C460").Value = "Maule"
$ws.Range("D460").Value = 44610
$ws.Range("E460").Value = 7
$ws.Range("F460").Value = "Fruta"
$ws.Range("G460").Value = 100101
$ws.Range("H460").Value = "Berries"
$ws.Range("I460").Value = 100112025
$ws.Range("J460").Value = "Frutilla"
$ws.Range("K460").Value = "Sin especificar"
$ws.Range("L460").Value = "Primera"
$ws.Range("M460").Value = 100
$ws.Range("N460").Value = 6000
$ws.Range("O460").Value = 6000
$ws.Range("P460").Value = 6000
$ws.Range("Q460").Value = "`$/bandeja 7 kilos"
$ws.Range("R460").Value = "Región del Maule"
$ws.Range("S460").Value = 857
$ws.Range("T460").Value = 7
